$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1897.8889
$ws.Range("I17").Value = 1169.5
$ws.Range("J17").Value = 1988.9375
$ws.Range("K17").Value = 3508.5
$ws.Range("L17").Value = 5966.8125
$ws.Range("M17").Value = -3340.5
$ws.Range("N17").Value = -6302.8125
$ws.Range("H131").Value = 4923.3
$ws.Range("I131").Value = 4923.3
$ws.Range("K131").Value = 14769.9
$ws.Range("M131").Value = -9729.900000000001
$ws.Range("H132").Value = 2052690.2
$ws.Range("I132").Value = 2317158
$ws.Range("K132").Value = 6951474
$ws.Range("M132").Value = -6948944
$ws.Range("H137").Value = 9340.209000000001
$ws.Range("I137").Value = 11387.064
$ws.Range("K137").Value = 34161.192
$ws.Range("M137").Value = -31611.192
$ws.Range("H138").Value = 1599.3334
$ws.Range("I138").Value = 1094.6744
$ws.Range("J138").Value = 3572.0908
$ws.Range("K138").Value = 3284.023200000001
$ws.Range("L138").Value = 10716.2724
$ws.Range("M138").Value = 1855.976799999999
$ws.Range("N138").Value = -20996.2724
$ws.Range("H141").Value = 1009.5517
$ws.Range("I141").Value = 770.4815
$ws.Range("J141").Value = 4237
$ws.Range("K141").Value = 2311.4445
$ws.Range("L141").Value = 12711
$ws.Range("M141").Value = 2868.5555
$ws.Range("N141").Value = -23071

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14082.325
$ws.Range("I32").Value = 14434.711
$ws.Range("K32").Value = 14434.711
$ws.Range("M32").Value = -14147.711
$ws.Range("H74").Value = 166046.78
$ws.Range("I74").Value = 215397.78
$ws.Range("J74").Value = 12510.333
$ws.Range("K74").Value = 215397.78
$ws.Range("L74").Value = 12510.333
$ws.Range("M74").Value = -214523.78
$ws.Range("N74").Value = -14258.333
$ws.Range("H77").Value = 166046.78
$ws.Range("I77").Value = 215397.78
$ws.Range("J77").Value = 12510.333
$ws.Range("K77").Value = 1076988.9
$ws.Range("L77").Value = 62551.665
$ws.Range("M77").Value = -1072620.9
$ws.Range("N77").Value = -71287.66500000001
$ws.Range("H102").Value = 7914.6055
$ws.Range("I102").Value = 8185.2905
$ws.Range("K102").Value = 8185.2905
$ws.Range("M102").Value = -6563.2905

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 74962.336
$ws.Range("J50").Value = 74962.336
$ws.Range("L50").Value = 74962.336
$ws.Range("N50").Value = -76110.336
$ws.Range("H99").Value = 2891.9143
$ws.Range("I99").Value = 2723
$ws.Range("K99").Value = 2723
$ws.Range("M99").Value = -1225
$ws.Range("H105").Value = 9616.733
$ws.Range("I105").Value = 9616.733
$ws.Range("K105").Value = 9616.733
$ws.Range("M105").Value = -7869.733

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 34999
$ws.Range("J51").Value = 34999
$ws.Range("L51").Value = 34999
$ws.Range("N51").Value = -36471
$ws.Range("H58").Value = 9049.048000000001
$ws.Range("I58").Value = 1225.3158
$ws.Range("K58").Value = 1225.3158
$ws.Range("M58").Value = -1022.3158
$ws.Range("H61").Value = 34999
$ws.Range("J61").Value = 34999
$ws.Range("L61").Value = 34999
$ws.Range("N61").Value = -35695
$ws.Range("H99").Value = 14121.625
$ws.Range("I99").Value = 10596.8
$ws.Range("J99").Value = 19996.334
$ws.Range("K99").Value = 10596.8
$ws.Range("L99").Value = 19996.334
$ws.Range("M99").Value = -9098.799999999999
$ws.Range("N99").Value = -22992.334
$ws.Range("H122").Value = 1737.875
$ws.Range("I122").Value = 1557.7142
$ws.Range("K122").Value = 4673.142599999999
$ws.Range("M122").Value = -2223.142599999999
$ws.Range("H126").Value = 14121.625
$ws.Range("I126").Value = 10596.8
$ws.Range("J126").Value = 19996.334
$ws.Range("K126").Value = 31790.4
$ws.Range("L126").Value = 59989.00199999999
$ws.Range("M126").Value = -29320.4
$ws.Range("N126").Value = -64929.00199999999
$ws.Range("H132").Value = 67859.87
$ws.Range("I132").Value = 111966.664
$ws.Range("K132").Value = 335899.992
$ws.Range("M132").Value = -333369.992
$ws.Range("H136").Value = 9049.048000000001
$ws.Range("I136").Value = 1225.3158
$ws.Range("K136").Value = 3675.9474
$ws.Range("M136").Value = -1125.9474
$ws.Range("H140").Value = 119949.89
$ws.Range("J140").Value = 119949.89
$ws.Range("L140").Value = 119949.89
$ws.Range("N140").Value = -130309.89

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 30129.5
$ws.Range("I109").Value = 30129.5
$ws.Range("K109").Value = 30129.5
$ws.Range("M109").Value = -29089.5
$ws.Range("H132").Value = 3543.1667
$ws.Range("I132").Value = 3639.8823
$ws.Range("K132").Value = 10919.6469
$ws.Range("M132").Value = -8389.6469

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2146.5789
$ws.Range("I7").Value = 1928.6471
$ws.Range("K7").Value = 1928.6471
$ws.Range("M7").Value = -1816.6471
$ws.Range("H40").Value = 2753
$ws.Range("I40").Value = 2725.5557
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2725.5557
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2589.5557
$ws.Range("N40").Value = -3272
$ws.Range("H68").Value = 3016.6086
$ws.Range("J68").Value = 4168
$ws.Range("L68").Value = 4168
$ws.Range("N68").Value = -5666
$ws.Range("H71").Value = 3016.6086
$ws.Range("J71").Value = 4168
$ws.Range("L71").Value = 20840
$ws.Range("N71").Value = -28328
$ws.Range("H93").Value = 5667.1665
$ws.Range("I93").Value = 3500
$ws.Range("K93").Value = 3500
$ws.Range("M93").Value = -2252
$ws.Range("H96").Value = 36666.668
$ws.Range("J96").Value = 36666.668
$ws.Range("L96").Value = 36666.668
$ws.Range("N96").Value = -42158.668
$ws.Range("H122").Value = 3314.3333
$ws.Range("I122").Value = 3537.8333
$ws.Range("J122").Value = 2867.3333
$ws.Range("K122").Value = 10613.4999
$ws.Range("L122").Value = 8601.999899999999
$ws.Range("M122").Value = -8163.499899999999
$ws.Range("N122").Value = -13501.9999
$ws.Range("H126").Value = 2146.5789
$ws.Range("I126").Value = 1928.6471
$ws.Range("K126").Value = 5785.9413
$ws.Range("M126").Value = -3315.9413
$ws.Range("H131").Value = 89999.5
$ws.Range("J131").Value = 89999.5
$ws.Range("L131").Value = 89999.5
$ws.Range("N131").Value = -100079.5
$ws.Range("H132").Value = 3376.7144
$ws.Range("I132").Value = 3522.4285
$ws.Range("K132").Value = 10567.2855
$ws.Range("M132").Value = -8037.2855
$ws.Range("H136").Value = 3239.4285
$ws.Range("I136").Value = 3302.1333
$ws.Range("J136").Value = 3082.6667
$ws.Range("K136").Value = 9906.3999
$ws.Range("L136").Value = 9248.000100000001
$ws.Range("M136").Value = -7356.3999
$ws.Range("N136").Value = -14348.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 28998.334
$ws.Range("J19").Value = 28998.334
$ws.Range("L19").Value = 28998.334
$ws.Range("N19").Value = -29346.334
$ws.Range("H100").Value = 1669.6666
$ws.Range("J100").Value = 2998
$ws.Range("L100").Value = 5996
$ws.Range("N100").Value = -7078
$ws.Range("H107").Value = 980.13336
$ws.Range("I107").Value = 636
$ws.Range("K107").Value = 1908
$ws.Range("M107").Value = 12
$ws.Range("H122").Value = 38993.555
$ws.Range("I122").Value = 44219.97
$ws.Range("J122").Value = 6589.8
$ws.Range("K122").Value = 132659.91
$ws.Range("L122").Value = 19769.4
$ws.Range("M122").Value = -130209.91
$ws.Range("N122").Value = -24669.4
$ws.Range("H126").Value = 229724.45
$ws.Range("I126").Value = 2089.75
$ws.Range("J126").Value = 836750.3
$ws.Range("K126").Value = 6269.25
$ws.Range("L126").Value = 2510250.9
$ws.Range("M126").Value = -3799.25
$ws.Range("N126").Value = -2515190.9
$ws.Range("H132").Value = 16862.04
$ws.Range("I132").Value = 23034.629
$ws.Range("K132").Value = 69103.887
$ws.Range("M132").Value = -66573.887
